# Adds four new "Atrações" (attractions) rows to the sheet:
#   - Advocate's Close
#   - Calton Hill
#   - St James Quarter
# and re-sorts "The Writers' Museum" / "Royal Mile" so the final order
# (rows 32-36) becomes:
#   32 The Writers' Museum
#   33 Advocate's Close      (new)
#   34 Royal Mile
#   35 Calton Hill           (new)
#   36 St James Quarter      (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert three blank rows right before the current row 32
# (Royal Mile). Each new row inherits formatting from row 31 above it,
# which already carries the "categoria"/"nome" column styles (A: s=11,
# B: s=6) used throughout this block of data.
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(32).Insert()

# Column G inherited an unwanted font style from row 31 ("Victoria
# Street"); the target rows use the plain default style there, so reset
# it back to Normal on the three freshly-inserted rows.
$ws.Cells.Item(32, 7).Style = "Normal"
$ws.Cells.Item(33, 7).Style = "Normal"
$ws.Cells.Item(34, 7).Style = "Normal"

# Row 32 will hold "The Writers' Museum", whose longitude cell uses a
# 6-decimal custom number format.
$ws.Cells.Item(32, 4).NumberFormat = "#,##0.000000"

# Row 36 currently still holds the old "The Writers' Museum" values
# (shifted down from row 33) including that same 6-decimal longitude
# format; it will be overwritten with "St James Quarter" data below, so
# drop the custom format back to the plain default style first.
$ws.Cells.Item(36, 4).Style = "Normal"

function Set-Row($r, $categoria, $nome, $lat, $lon, $icone, $endereco, $descricao) {
    $ws.Cells.Item($r, 1).Value = $categoria
    $ws.Cells.Item($r, 2).Value = $nome
    $ws.Cells.Item($r, 3).Value = $lat
    $ws.Cells.Item($r, 4).Value = $lon
    $ws.Cells.Item($r, 5).Value = $icone
    $ws.Cells.Item($r, 6).Value = $endereco
    $ws.Cells.Item($r, 7).Value = $descricao
}

# Same as Set-Row, but writes the "descricao" (G) column before the
# "endereco" (F) column - matches the order the new shared strings were
# authored in for this row.
function Set-Row-DescFirst($r, $categoria, $nome, $lat, $lon, $icone, $endereco, $descricao) {
    $ws.Cells.Item($r, 1).Value = $categoria
    $ws.Cells.Item($r, 2).Value = $nome
    $ws.Cells.Item($r, 3).Value = $lat
    $ws.Cells.Item($r, 4).Value = $lon
    $ws.Cells.Item($r, 5).Value = $icone
    $ws.Cells.Item($r, 7).Value = $descricao
    $ws.Cells.Item($r, 6).Value = $endereco
}

Set-Row 32 "Atrações" "The Writers' Museum" 55.949809433784303 -3.1936380192635001 `
    "touristic.png" "Lawnmarket, Lady Stair's Cl, Edinburgh EH1 2PA, Reino Unido" `
    "Manuscritos e itens pessoais de Burns, Scott e Stevenson, exibidos em uma casa do século 17."

Set-Row 33 "Atrações" "Advocate's Close" 55.949785395700403 -3.1913097995194399 `
    "touristic.png" "WRX5+WF Edimburgo, Reino Unido" `
    "Advocates Close é um beco estreito e íngreme em Edimburgo, de origem medieval, reconstruído no início do século XXI."

Set-Row 34 "Atrações" "Royal Mile" 55.949604663396798 -3.1917009894855899 `
    "touristic.png" "Edinburgh EH1 1QS, Reino Unido" `
    "rua principal que liga o Castelo de Edimburgo ao Palácio de Holyrood. Muitas lojinhas de tartan, pubs e gaitas de fole."

Set-Row-DescFirst 35 "Atrações" "Calton Hill " 55.955621712536299 -3.1821399240250501 `
    "touristic.png" "Edimburgo EH7 5AA, Reino Unido" `
    "Colina proeminente com monumentos neoclássicos e vista panorâmica para o horizonte da cidade e o Castelo de Edimburgo."

Set-Row-DescFirst 36 "Atrações" "St James Quarter " 55.955367742394401 -3.1885290367581001 `
    "touristic.png" "St James Cres, Edinburgh EH1 3AD, Reino Unido" `
    "Shopping contemporâneo com grandes marcas de varejo, diversos restaurantes, bares e cafés."

# Match the author's final selection / scroll position.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A36").Select()
